# Update dSF (column F) values on Sheet1 to reflect the repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -6
$ws.Range("F4").Value = -7
$ws.Range("F6").Value = -14
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -4
